# Apply "update building block types" changes to the GEO High-throughput sequencing Assay template

$wb = $excel.ActiveWorkbook

# --- isa_template sheet (metadata) ---
$ws1 = $wb.Worksheets.Item("isa_template")

# Version: 1.0.2 -> 1.0.3
$ws1.Range("B4").Value2 = "1.0.3"

# Table name: Assay -> New Table
$ws1.Range("B7").Value2 = "New Table"

# ER Term Accession Number: obolibrary purl -> nfdi4plants purl
$ws1.Range("B10").Value2 = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_1000096"

# Tags Term Accession Number: obolibrary purls -> bioregistry.io
$ws1.Range("B14").Value2 = "https://bioregistry.io/EFO:0004184"
$ws1.Range("C14").Value2 = "https://bioregistry.io/EFO:0008896"
$ws1.Range("D14").Value2 = "https://bioregistry.io/NCIT:C153189"

# --- Assay (table) sheet ---
$ws2 = $wb.Worksheets.Item("Assay")

# Rename the sheet itself: Assay -> New Table
$ws2.Name = "New Table"

# Header row: rename building block columns
$ws2.Range("M1").Value2 = "Term Source REF (GENEPIO:0001973)"
$ws2.Range("N1").Value2 = "Term Accession Number (GENEPIO:0001973)"
$ws2.Range("O1").Value2 = "Component [next generation sequencing instrument model]"
$ws2.Range("R1").Value2 = "Output [Data]"

# Data row (row 2): update values
$ws2.Range("C2").Value2 = ""
$ws2.Range("K2").Value2 = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_0000086"
$ws2.Range("N2").Value2 = "https://bioregistry.io/EFO:0008896"
$ws2.Range("Q2").Value2 = "https://bioregistry.io/OBI:0002001"

$wb.Save()
